# Weekly update: insert a new record at row 10, shifting existing rows 10-69 down to 11-70.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 10 (shifts rows 10:69 down to 11:70, copies formatting from row above)
$ws.Rows.Item(10).Insert()

# Populate the new row 10 with this week's data
$ws.Cells.Item(10, 1).Value  = 2
$ws.Cells.Item(10, 2).Value  = "Comercializadora del Agro de Limarí"
$ws.Cells.Item(10, 3).Value  = "Coquimbo"
$ws.Cells.Item(10, 4).Value  = 44959
$ws.Cells.Item(10, 5).Value  = 4
$ws.Cells.Item(10, 6).Value  = 100112032
$ws.Cells.Item(10, 7).Value  = "Zapallo italiano"
$ws.Cells.Item(10, 8).Value  = "Sin especificar"
$ws.Cells.Item(10, 9).Value  = "Primera"
$ws.Cells.Item(10, 10).Value = 500
$ws.Cells.Item(10, 11).Value = 6000
$ws.Cells.Item(10, 12).Value = 7000
$ws.Cells.Item(10, 13).Value = 6500
$ws.Cells.Item(10, 14).Value = "$/caja 70 unidades"
$ws.Cells.Item(10, 15).Value = "Provincia de Limarí"
$ws.Cells.Item(10, 16).Value = 93
$ws.Cells.Item(10, 17).Value = 70
$ws.Cells.Item(10, 18).Value = "Hortaliza"
